# "running profile test cases"
#
# On the "Test Cases" sheet, the ProfileFollowerTest row (row 12) has its
# Runmode flipped from "N" to "Y" so that test case is actually executed
# in the run. Leave the selection where the author left it (B15).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$ws.Range("C12").Value = "Y"

$ws.Range("B15").Select() | Out-Null
